# The deck's design/theme was switched from the custom "Integral" theme to
# the default "Office Theme". In the underlying OOXML this shows up as the
# active theme part (ppt/theme/theme1.xml, the one wired to the slide
# master) getting the Office Theme's 12-color scheme (clrScheme) - the
# font scheme and format scheme are already identical between the two
# themes in this file, so only the colors need to change.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

function Set-ThemeColor {
    param([int]$Index, [int]$R, [int]$G, [int]$B)
    $tcs.Colors($Index).RGB = $R + ($G * 256) + ($B * 65536)
}

# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
Set-ThemeColor 1  0x00 0x00 0x00   # dk1
Set-ThemeColor 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor 12 0x95 0x4F 0x72   # folHlink
